# 01: Finished experimental part
# Adds the "relative error" propagation columns (R/S and Z/AA) for rows
# 30/32/34, the combined-error row 34 (P/Q and X/Y), the "U Ar [eV]"
# literature-value comparison in row 36, and the row-37 note, on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 30: relative error of the first exp-fit error term -------------
$ws.Range("R30").Value = "relative error"
$ws.Range("S30").Formula = "=Q30/Q29"

$ws.Range("Z30").Value = "relative error"
$ws.Range("AA30").Formula = "=Y30/Y29"

# --- Row 32: relative error of the second exp-fit error term ------------
$ws.Range("R32").Value = "relative error"
$ws.Range("S32").Formula = "=Q32/Q31"

$ws.Range("Z32").Value = "relative error"
$ws.Range("AA32").Formula = "=Y32/Y31"

# --- Row 34: combined (propagated) error of alpha/p ----------------------
$ws.Range("P34").Value = "error"
$ws.Range("P34").NumberFormat = "0.000"
$ws.Range("Q34").Formula = "=Q33*S34"

$ws.Range("R34").Value = "relative error"
$ws.Range("S34").Formula = "=SQRT((S30^2)+(S32^2))"

$ws.Range("X34").Value = "error"
$ws.Range("X34").NumberFormat = "0.000"
$ws.Range("Y34").Formula = "=Y33*AA34"

$ws.Range("Z34").Value = "relative error"
$ws.Range("AA34").Formula = "=SQRT((AA30^2)+(AA32^2))"

# --- Row 36: comparison against the Ar ionisation literature value ------
$ws.Range("N36").Value = "U Ar [eV]"
$ws.Range("O36").Value = 15.76
$ws.Range("Q36").Formula = "=O36-Q33"
$ws.Range("Y36").Formula = "=Y33-O36"

# --- Row 37: note ---------------------------------------------------------
$ws.Range("Q37").Value = "přesnější, viz error"

# --- Column width for the new "relative error" label column (R) ---------
$ws.Columns.Item(18).ColumnWidth = 12.0221354166667

# --- View state: scrolled/selected position after the edits -------------
$null = $ws.Range("O36").Select()
$excel.ActiveWindow.ScrollColumn = 11

$wb.Saved = $false
